$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I3").Value = "11g"
$ws.Range("A3:J3").HorizontalAlignment = -4108
$ws.Range("H14").Select() | Out-Null
